# Update automàtic: dades i banners [2026-02-24 23:50]
# Refresh the per-station extraction timestamps and the associated
# re-measured readings (humidity %, radiation, pressure, temperatures)
# coming from the new meteo.cat poll.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("E2").Value = "2026-02-24 23:48:42"
$ws.Range("H2").NumberFormat = "@"
$ws.Range("H2").Value = "43%"
$ws.Range("K2").Value = "12.6 MJ/m2"
$ws.Range("O2").Value = "5.5 °C"
# Row 3
$ws.Range("E3").Value = "2026-02-24 23:48:44"
$ws.Range("H3").NumberFormat = "@"
$ws.Range("H3").Value = "34%"
$ws.Range("K3").Value = "16.3 MJ/m2"
# Row 4
$ws.Range("E4").Value = "2026-02-24 23:48:47"
$ws.Range("H4").NumberFormat = "@"
$ws.Range("H4").Value = "75%"
$ws.Range("O4").Value = "12.1 °C"
# Row 5
$ws.Range("E5").Value = "2026-02-24 23:48:49"
$ws.Range("O5").Value = "5.7 °C"
# Row 6
$ws.Range("E6").Value = "2026-02-24 23:48:52"
$ws.Range("H6").NumberFormat = "@"
$ws.Range("H6").Value = "73%"
$ws.Range("O6").Value = "13.7 °C"
# Row 7
$ws.Range("E7").Value = "2026-02-24 23:48:54"
$ws.Range("H7").NumberFormat = "@"
$ws.Range("H7").Value = "75%"
# Row 8
$ws.Range("E8").Value = "2026-02-24 23:48:56"
# Row 9
$ws.Range("E9").Value = "2026-02-24 23:48:58"
$ws.Range("O9").Value = "11.4 °C"
# Row 10
$ws.Range("E10").Value = "2026-02-24 23:49:01"
$ws.Range("O10").Value = "10.6 °C"
# Row 11
$ws.Range("E11").Value = "2026-02-24 23:49:03"
# Row 12
$ws.Range("E12").Value = "2026-02-24 23:49:05"
$ws.Range("O12").Value = "10.2 °C"
# Row 13
$ws.Range("E13").Value = "2026-02-24 23:49:07"
$ws.Range("O13").Value = "6.4 °C"
# Row 14
$ws.Range("E14").Value = "2026-02-24 23:49:10"
$ws.Range("N14").Value = "4.2 °C 23:29 TU"
$ws.Range("O14").Value = "10.9 °C"
# Row 15
$ws.Range("E15").Value = "2026-02-24 23:49:12"
$ws.Range("O15").Value = "11.4 °C"
# Row 16
$ws.Range("E16").Value = "2026-02-24 23:49:14"
# Row 17
$ws.Range("E17").Value = "2026-02-24 23:49:17"
# Row 18
$ws.Range("E18").Value = "2026-02-24 23:49:20"
$ws.Range("O18").Value = "10.8 °C"
# Row 19
$ws.Range("E19").Value = "2026-02-24 23:49:23"
# Row 20
$ws.Range("E20").Value = "2026-02-24 23:49:25"
$ws.Range("O20").Value = "3.5 °C"
# Row 21
$ws.Range("E21").Value = "2026-02-24 23:49:27"
$ws.Range("O21").Value = "9.5 °C"
# Row 22
$ws.Range("E22").Value = "2026-02-24 23:49:30"
$ws.Range("H22").NumberFormat = "@"
$ws.Range("H22").Value = "28%"
$ws.Range("K22").Value = "16.3 MJ/m2"
$ws.Range("N22").Value = "0.4 °C 23:26 TU"
$ws.Range("O22").Value = "3.3 °C"
# Row 23
$ws.Range("E23").Value = "2026-02-24 23:49:33"
$ws.Range("O23").Value = "4.6 °C"
# Row 24
$ws.Range("E24").Value = "2026-02-24 23:49:35"
$ws.Range("J24").Value = "1020.9 hPa"
# Row 25
$ws.Range("E25").Value = "2026-02-24 23:49:38"
# Row 26
$ws.Range("E26").Value = "2026-02-24 23:49:41"
$ws.Range("H26").NumberFormat = "@"
$ws.Range("H26").Value = "42%"
# Row 27
$ws.Range("E27").Value = "2026-02-24 23:49:44"
$ws.Range("H27").NumberFormat = "@"
$ws.Range("H27").Value = "34%"
$ws.Range("O27").Value = "6.1 °C"
# Row 28
$ws.Range("E28").Value = "2026-02-24 23:49:46"
$ws.Range("O28").Value = "11.1 °C"
# Row 29
$ws.Range("E29").Value = "2026-02-24 23:49:49"
$ws.Range("O29").Value = "9.9 °C"
# Row 30
$ws.Range("E30").Value = "2026-02-24 23:49:52"
$ws.Range("H30").NumberFormat = "@"
$ws.Range("H30").Value = "78%"
$ws.Range("O30").Value = "12.7 °C"
# Row 31
$ws.Range("E31").Value = "2026-02-24 23:49:55"
$ws.Range("J31").Value = "1018.9 hPa"
$ws.Range("O31").Value = "15.5 °C"
# Row 32
$ws.Range("E32").Value = "2026-02-24 23:49:57"
$ws.Range("H32").NumberFormat = "@"
$ws.Range("H32").Value = "71%"
$ws.Range("O32").Value = "6.6 °C"
# Row 33
$ws.Range("E33").Value = "2026-02-24 23:50:00"
$ws.Range("J33").Value = "1021.7 hPa"
$ws.Range("O33").Value = "8.3 °C"
# Row 34
$ws.Range("E34").Value = "2026-02-24 23:50:03"
$ws.Range("O34").Value = "4.4 °C"
# Row 35
$ws.Range("E35").Value = "2026-02-24 23:50:05"
$ws.Range("K35").Value = "16.6 MJ/m2"
# Row 36
$ws.Range("E36").Value = "2026-02-24 23:50:08"
$ws.Range("H36").NumberFormat = "@"
$ws.Range("H36").Value = "83%"
$ws.Range("O36").Value = "12.6 °C"
# Row 37
$ws.Range("E37").Value = "2026-02-24 23:50:11"
$ws.Range("H37").NumberFormat = "@"
$ws.Range("H37").Value = "74%"
$ws.Range("O37").Value = "8.2 °C"
# Row 38
$ws.Range("E38").Value = "2026-02-24 23:50:13"
$ws.Range("K38").Value = "14.5 MJ/m2"
$ws.Range("O38").Value = "11.5 °C"
# Row 39
$ws.Range("E39").Value = "2026-02-24 23:50:16"
# Row 40
$ws.Range("E40").Value = "2026-02-24 23:50:18"
$ws.Range("H40").NumberFormat = "@"
$ws.Range("H40").Value = "68%"
$ws.Range("O40").Value = "8.1 °C"
# Row 41
$ws.Range("E41").Value = "2026-02-24 23:50:21"
$ws.Range("H41").NumberFormat = "@"
$ws.Range("H41").Value = "82%"
$ws.Range("J41").Value = "1020.3 hPa"
# Row 42
$ws.Range("E42").Value = "2026-02-24 23:50:23"
$ws.Range("O42").Value = "11.0 °C"
# Row 43
$ws.Range("E43").Value = "2026-02-24 23:50:25"
$ws.Range("H43").NumberFormat = "@"
$ws.Range("H43").Value = "70%"
$ws.Range("O43").Value = "10.4 °C"
# Row 44
$ws.Range("E44").Value = "2026-02-24 23:50:28"
$ws.Range("O44").Value = "2.7 °C"
# Row 45
$ws.Range("E45").Value = "2026-02-24 23:50:31"
# Row 46
$ws.Range("E46").Value = "2026-02-24 23:50:34"
$ws.Range("K46").Value = "15.5 MJ/m2"
$ws.Range("O46").Value = "10.4 °C"
